$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# Insert a new column before column L ("Biological risk assessment number"),
# shifting the existing HuMFre..Section thickness columns one to the right.
$ws.Columns("L:L").Insert()

# Populate the new column. Values are entered in this order so the shared
# string table is built up in the same sequence as the source workbook.
$ws.Range("L2").Value = "Biological risk assessment number"
$ws.Range("L4").Value = "RISK1"
$ws.Range("L5").Value = "RISK3"
$ws.Range("L3").Value = "RISKX"

# Copy the formatting from the (now shifted) neighbouring column so the new
# column's cells carry the same style as the rest of their row.
$ws.Range("M2:M5").Copy()
$ws.Range("L2:L5").PasteSpecial(-4122)

# Update the active selection to match the saved workbook state.
$null = $ws.Range("L6").Select()
